$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.301.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "'1.872.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "'1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'0.7098"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'241.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +2.02%  "

$ws.Range("D9").Value = "'0.3098"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").Value = "'25.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("D11").Value = "'0.08406"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

$ws.Range("D12").Value = "'1.871.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").Value = "'5.243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "'0.7112"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").Value = "'91.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "'29.312.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'6.074"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.11%  "

$ws.Range("D18").Value = "'0.000008197"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.71%  "

$ws.Range("D19").Value = "'240.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "

$ws.Range("D20").Value = "'13.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("D21").Value = "'2.124.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'7.764"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'0.1592"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("D26").Value = "'163.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("D28").Value = "'18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "'1.504"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'4.394"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "'4.296"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("D33").Value = "'0.05377"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.37%  "

$ws.Range("D34").Value = "'1.942"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("D35").Value = "'1.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("D36").Value = "'0.7491"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.85%  "

$ws.Range("D37").Value = "'2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "

$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("D39").Value = "'1.230.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.07%  "

$ws.Range("D40").Value = "'2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "

$ws.Range("D41").Value = "'6.548"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.51%  "

$ws.Range("D42").Value = "'0.8871"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("D43").Value = "'72.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("D44").Value = "'108.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.00%  "

$ws.Range("D45").Value = "'1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "'2.020.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").Value = "'0.5193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").Value = "'1.792"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("E49").Value = "  +2.10%  "

$ws.Range("D50").Value = "'9.410"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.80%  "

$ws.Range("D51").Value = "'0.4313"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
